$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row2
$ws.Range("G2").Value = 3.157093666666667
$ws.Range("H2").Value = 9.471281000000001
$ws.Range("I2").Value = 0.8981781966433163
$ws.Range("J2").Value = 0.8981781966433162
$ws.Range("M2").Value = 14.18473066666667
$ws.Range("N2").Value = 42.554192
$ws.Range("O2").Value = 0.2161726484729168
$ws.Range("P2").Value = 0.2161726484729169
$ws.Range("Q2").Value = 44.78252335110578
$ws.Range("R2").Value = 403.042710159952
$ws.Range("S2").Value = 0.194161559569014
$ws.Range("T2").Value = 0.194161559569014

# Row3
$ws.Range("G3").Value = 3.157093666666667
$ws.Range("H3").Value = 9.471281000000001
$ws.Range("I3").Value = 0.8981781966433163
$ws.Range("J3").Value = 0.8981781966433162
$ws.Range("O3").Value = 0.3233301721219069
$ws.Range("P3").Value = 0.3233301721219069
$ws.Range("Q3").Value = 66.98137384841456
$ws.Range("R3").Value = 602.8323646357311
$ws.Range("S3").Value = 0.2904081109168274
$ws.Range("T3").Value = 0.2904081109168274

# Row4
$ws.Range("G4").Value = 3.157093666666667
$ws.Range("H4").Value = 9.471281000000001
$ws.Range("I4").Value = 0.8981781966433163
$ws.Range("J4").Value = 0.8981781966433162
$ws.Range("M4").Value = 15.72529
$ws.Range("N4").Value = 47.17587
$ws.Range("O4").Value = 0.2396504852427705
$ws.Range("P4").Value = 0.2396504852427705
$ws.Range("Q4").Value = 49.64621346549667
$ws.Range("R4").Value = 446.81592118947
$ws.Range("S4").Value = 0.2152488406600473
$ws.Range("T4").Value = 0.2152488406600473

# Row5
$ws.Range("G5").Value = 3.157093666666667
$ws.Range("H5").Value = 9.471281000000001
$ws.Range("I5").Value = 0.8981781966433163
$ws.Range("J5").Value = 0.8981781966433162
$ws.Range("M5").Value = 4.273701999999999
$ws.Range("N5").Value = 12.821106
$ws.Range("O5").Value = 0.0651304210022835
$ws.Range("P5").Value = 0.06513042100228351
$ws.Range("Q5").Value = 13.49247751742067
$ws.Range("R5").Value = 121.432297656786
$ws.Range("S5").Value = 0.05849872408245097
$ws.Range("T5").Value = 0.05849872408245097

# Row6
$ws.Range("G6").Value = 3.157093666666667
$ws.Range("H6").Value = 9.471281000000001
$ws.Range("I6").Value = 0.8981781966433163
$ws.Range("J6").Value = 0.8981781966433162
$ws.Range("M6").Value = 10.21772833333333
$ws.Range("N6").Value = 30.653185
$ws.Range("O6").Value = 0.1557162731601222
$ws.Range("P6").Value = 0.1557162731601222
$ws.Range("Q6").Value = 32.25832540888722
$ws.Range("R6").Value = 290.324928679985
$ws.Range("S6").Value = 0.1398609614149766
$ws.Range("T6").Value = 0.1398609614149766

# Row7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.3579033333333333
$ws.Range("H7").Value = 1.07371
$ws.Range("I7").Value = 0.1018218033566837
$ws.Range("J7").Value = 0.1018218033566837
$ws.Range("M7").Value = 14.18473066666667
$ws.Range("N7").Value = 42.554192
$ws.Range("O7").Value = 0.2161726484729168
$ws.Range("P7").Value = 0.2161726484729169
$ws.Range("Q7").Value = 5.076762388035555
$ws.Range("R7").Value = 45.69086149232
$ws.Range("S7").Value = 0.02201108890390286
$ws.Range("T7").Value = 0.02201108890390286

# Row8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.3579033333333333
$ws.Range("H8").Value = 1.07371
$ws.Range("I8").Value = 0.1018218033566837
$ws.Range("J8").Value = 0.1018218033566837
$ws.Range("O8").Value = 0.3233301721219069
$ws.Range("P8").Value = 0.3233301721219069
$ws.Range("Q8").Value = 7.593330924801111
$ws.Range("R8").Value = 68.33997832321
$ws.Range("S8").Value = 0.03292206120507951
$ws.Range("T8").Value = 0.03292206120507951

# Row9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.3579033333333333
$ws.Range("H9").Value = 1.07371
$ws.Range("I9").Value = 0.1018218033566837
$ws.Range("J9").Value = 0.1018218033566837
$ws.Range("M9").Value = 15.72529
$ws.Range("N9").Value = 47.17587
$ws.Range("O9").Value = 0.2396504852427705
$ws.Range("P9").Value = 0.2396504852427705
$ws.Range("Q9").Value = 5.628133708633333
$ws.Range("R9").Value = 50.6532033777
$ws.Range("S9").Value = 0.02440164458272322
$ws.Range("T9").Value = 0.02440164458272322

# Row10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.3579033333333333
$ws.Range("H10").Value = 1.07371
$ws.Range("I10").Value = 0.1018218033566837
$ws.Range("J10").Value = 0.1018218033566837
$ws.Range("M10").Value = 4.273701999999999
$ws.Range("N10").Value = 12.821106
$ws.Range("O10").Value = 0.0651304210022835
$ws.Range("P10").Value = 0.06513042100228351
$ws.Range("Q10").Value = 1.529572191473333
$ws.Range("R10").Value = 13.76614972326
$ws.Range("S10").Value = 0.006631696919832535
$ws.Range("T10").Value = 0.006631696919832537

# Row11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.3579033333333333
$ws.Range("H11").Value = 1.07371
$ws.Range("I11").Value = 0.1018218033566837
$ws.Range("J11").Value = 0.1018218033566837
$ws.Range("M11").Value = 10.21772833333333
$ws.Range("N11").Value = 30.653185
$ws.Range("O11").Value = 0.1557162731601222
$ws.Range("P11").Value = 0.1557162731601222
$ws.Range("Q11").Value = 3.656959029594444
$ws.Range("R11").Value = 32.91263126635
$ws.Range("S11").Value = 0.01585531174514561
$ws.Range("T11").Value = 0.01585531174514562
